# Renumber the parish ("paroisse") group codes in column C from the old
# 0305020XXX scheme to the new 0403110XXX scheme (constant +98090000 offset,
# rows 2-23) and drop the now-unneeded explicit number/alignment format on
# that column so it falls back to the sheet's default style. Matches the
# commit "Corrige les numeros de groupes pour passer de 030502XXXX a
# 040311XXXX pour les paroisses".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$offset = 98090000

for ($row = 2; $row -le 23; $row++) {
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = $cCell.Value2 + $offset
    $cCell.ClearFormats()
}

# The active selection moved from G12 to D25.
$ws.Range("D25").Select() | Out-Null
